# Case_2_4 (380 kV) power-flow re-run: slack-bus voltage setpoint (B2:B25)
# dropped from 1.05 pu to 1.02 pu, which shifts every other bus voltage
# magnitude in the results table (columns C:F, I:N). Column G is the
# external-grid bus and stays at 1.0 pu; column H has no data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.011700888260167
$ws.Range("D2").Value = 1.018373784019337
$ws.Range("E2").Value = 0.9926147277508489
$ws.Range("F2").Value = 1.009963350290567
$ws.Range("I2").Value = 1.024737108760112
$ws.Range("J2").Value = 1.016947847641629
$ws.Range("K2").Value = 1.02122019089174
$ws.Range("L2").Value = 0.9955398523336033
$ws.Range("M2").Value = 1.012834977419816
$ws.Range("N2").Value = 1.009814735657774

# row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.012567919567097
$ws.Range("D3").Value = 1.018984946931742
$ws.Range("E3").Value = 0.9936372048519304
$ws.Range("F3").Value = 1.011475701136261
$ws.Range("I3").Value = 1.024848710347437
$ws.Range("J3").Value = 1.017448690233878
$ws.Range("K3").Value = 1.021637686750333
$ws.Range("L3").Value = 0.9963617723202692
$ws.Range("M3").Value = 1.014149312728175
$ws.Range("N3").Value = 1.009979454800132

# row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.01312849593714
$ws.Range("D4").Value = 1.019379526142619
$ws.Range("E4").Value = 0.9942998659930995
$ws.Range("F4").Value = 1.012453752266924
$ws.Range("I4").Value = 1.024918697281743
$ws.Range("J4").Value = 1.017771719002855
$ws.Range("K4").Value = 1.021906261951088
$ws.Range("L4").Value = 0.9968940712668345
$ws.Range("M4").Value = 1.014998758648693
$ws.Range("N4").Value = 1.010085685973224

# row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.013364054119463
$ws.Range("D5").Value = 1.019545194602445
$ws.Range("E5").Value = 0.9945786998346017
$ws.Range("F5").Value = 1.012864800561566
$ws.Range("I5").Value = 1.024947586105939
$ws.Range("J5").Value = 1.017907268768451
$ws.Range("K5").Value = 1.022018793876521
$ws.Range("L5").Value = 0.997117960005301
$ws.Range("M5").Value = 1.015355626569939
$ws.Range("N5").Value = 1.010130260942786

# row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.013403599028691
$ws.Range("D6").Value = 1.01957299856721
$ws.Range("E6").Value = 0.9946255319796338
$ws.Range("F6").Value = 1.012933810232236
$ws.Range("I6").Value = 1.024952405353728
$ws.Range("J6").Value = 1.01793001341511
$ws.Range("K6").Value = 1.02203766635093
$ws.Range("L6").Value = 0.9971555583673453
$ws.Range("M6").Value = 1.015415532365469
$ws.Range("N6").Value = 1.010137740312107

# row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.013131643902216
$ws.Range("D7").Value = 1.019381740649512
$ws.Range("E7").Value = 0.9943035907982488
$ws.Range("F7").Value = 1.012459245196837
$ws.Range("I7").Value = 1.024919085393408
$ws.Range("J7").Value = 1.017773531214059
$ws.Range("K7").Value = 1.021907767090704
$ws.Range("L7").Value = 0.9968970624462087
$ws.Range("M7").Value = 1.015003528069676
$ws.Range("N7").Value = 1.010086281918737

# row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.011993998796304
$ws.Range("D8").Value = 1.018580511845594
$ws.Range("E8").Value = 0.9929600610674301
$ws.Range("F8").Value = 1.01047457144109
$ws.Range("I8").Value = 1.0247752857505
$ws.Range("J8").Value = 1.017117327017547
$ws.Range("K8").Value = 1.021361610808137
$ws.Range("L8").Value = 0.995817528259106
$ws.Range("M8").Value = 1.013279376347651
$ws.Range("N8").Value = 1.00987047633582

# row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.009985876688987
$ws.Range("D9").Value = 1.017161916230851
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.006972949548811
$ws.Range("I9").Value = 1.024504863592926
$ws.Range("J9").Value = 1.015952977712818
$ws.Range("K9").Value = 1.02038719622062
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.010233218705661
$ws.Range("N9").Value = 1.009487497825462

# row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.008644817217841
$ws.Range("D10").Value = 1.016211717238475
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.004635271100022
$ws.Range("I10").Value = 1.024313165759963
$ws.Range("J10").Value = 1.015171350168007
$ws.Range("K10").Value = 1.019729539567003
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.008196806991192
$ws.Range("N10").Value = 1.009230365610914

# row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.008063575365778
$ws.Range("D11").Value = 1.015799221475528
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.003622183686665
$ws.Range("I11").Value = 1.02422745709784
$ws.Range("J11").Value = 1.014831617751766
$ws.Range("K11").Value = 1.019442864483395
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.007313618953825
$ws.Range("N11").Value = 1.009118594811778

# row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.007847592849912
$ws.Range("D12").Value = 1.015645844712157
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.003245742774031
$ws.Range("I12").Value = 1.024195215552719
$ws.Range("J12").Value = 1.014705233251505
$ws.Range("K12").Value = 1.019336094724189
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.006985346606095
$ws.Range("N12").Value = 1.00907701342331

# row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.007893925609761
$ws.Range("D13").Value = 1.015678751662113
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.00332649684541
$ws.Range("I13").Value = 1.024202149817689
$ws.Range("J13").Value = 1.014732351879627
$ws.Range("K13").Value = 1.019359010129944
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.007055772045235
$ws.Range("N13").Value = 1.009085935701928

# row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.008045723886103
$ws.Range("D14").Value = 1.015786546512336
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.003591069757086
$ws.Range("I14").Value = 1.024224800271838
$ws.Range("J14").Value = 1.014821174702127
$ws.Range("K14").Value = 1.019434044691194
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.00728648831699
$ws.Range("N14").Value = 1.009115159001906

# row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.008139240756685
$ws.Range("D15").Value = 1.015852941633835
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.00375406376197
$ws.Range("I15").Value = 1.024238702240504
$ws.Range("J15").Value = 1.014875875830476
$ws.Range("K15").Value = 1.019480238092022
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.007428611311582
$ws.Range("N15").Value = 1.009133155862951

# row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.008683380636064
$ws.Range("D16").Value = 1.016239071098364
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.004702487773215
$ws.Range("I16").Value = 1.024318797067157
$ws.Range("J16").Value = 1.015193870062629
$ws.Range("K16").Value = 1.019748525137655
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.00825539114045
$ws.Range("N16").Value = 1.009237774392741

# row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.009024556722395
$ws.Range("D17").Value = 1.016480998461306
$ws.Range("E17").Value = 0.989476357848556
$ws.Range("F17").Value = 1.005297175126365
$ws.Range("I17").Value = 1.024368315291226
$ws.Range("J17").Value = 1.015392996029741
$ws.Range("K17").Value = 1.019916304666739
$ws.Range("L17").Value = 0.9930127773699352
$ws.Range("M17").Value = 1.008773627434139
$ws.Range("N17").Value = 1.009303283478062

# row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.00922350549452
$ws.Range("D18").Value = 1.016622008916878
$ws.Range("E18").Value = 0.9897087662937556
$ws.Range("F18").Value = 1.005643963917157
$ws.Range("I18").Value = 1.024396937677607
$ws.Range("J18").Value = 1.015509019097165
$ws.Range("K18").Value = 1.020013983757049
$ws.Range("L18").Value = 0.9932001317071769
$ws.Range("M18").Value = 1.009075770377467
$ws.Range("N18").Value = 1.009341452248794

# row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.009291332819707
$ws.Range("D19").Value = 1.016670072567096
$ws.Range("E19").Value = 0.9897880325774034
$ws.Range("F19").Value = 1.00576219616936
$ws.Range("I19").Value = 1.024406652929138
$ws.Range("J19").Value = 1.015548558978214
$ws.Range("K19").Value = 1.020047258602617
$ws.Range("L19").Value = 0.9932640239640975
$ws.Range("M19").Value = 1.009178770490507
$ws.Range("N19").Value = 1.009354459762295

# row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.008987957288347
$ws.Range("D20").Value = 1.016455052456525
$ws.Range("E20").Value = 0.9894336180360679
$ws.Range("F20").Value = 1.005233379357276
$ws.Range("I20").Value = 1.024363029422307
$ws.Range("J20").Value = 1.015371644497355
$ws.Range("K20").Value = 1.019898322526685
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.008718039640073
$ws.Range("N20").Value = 1.009296259272996

# row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.008001025380997
$ws.Range("D21").Value = 1.015754807956456
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.003513163410628
$ws.Range("I21").Value = 1.02421814146876
$ws.Range("J21").Value = 1.014795023916337
$ws.Range("K21").Value = 1.019411956771367
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.007218554158861
$ws.Range("N21").Value = 1.009106555255569

# row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.007380019440693
$ws.Range("D22").Value = 1.015313625741673
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.002430810319801
$ws.Range("I22").Value = 1.024124698592994
$ws.Range("J22").Value = 1.01443136440953
$ws.Range("K22").Value = 1.019104505552012
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.006274510004032
$ws.Range("N22").Value = 1.008986906255054

# row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.007709272243503
$ws.Range("D23").Value = 1.015547590850355
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.003004662717958
$ws.Range("I23").Value = 1.024174456663089
$ws.Range("J23").Value = 1.014624252935517
$ws.Range("K23").Value = 1.019267647890193
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.00677508674112
$ws.Range("N23").Value = 1.009050369960999

# row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.009004495162197
$ws.Range("D24").Value = 1.016466776653412
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.005262206171984
$ws.Range("I24").Value = 1.02436541868462
$ws.Range("J24").Value = 1.015381292717962
$ws.Range("K24").Value = 1.0199064484498
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.008743157788867
$ws.Range("N24").Value = 1.009299433337311

# row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.010505432134435
$ws.Range("D25").Value = 1.017529448769823
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.007878754181918
$ws.Range("I25").Value = 1.024576788776056
$ws.Range("J25").Value = 1.016254941199935
$ws.Range("K25").Value = 1.02064052641772
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.01102169741244
$ws.Range("N25").Value = 1.009586826750631

